$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New SUM formula summarising column C, placed at I4 (cached value 120)
$ws.Range("I4").Formula = "=SUM(C:C)"

# Row 39: answer corrected to "red.png" and marked correct (D/E = 1)
$ws.Range("B39").Value = "red.png"
$ws.Range("D39").Value = 1
$ws.Range("E39").Value = 1

# Row 41: answer corrected to "orange.png" and duration fixed to 3
$ws.Range("B41").Value = "orange.png"
$ws.Range("C41").Value = 3

# New trailing row 42 appended to the data table
$ws.Range("A42").Value = 1
$ws.Range("B42").Value = "green.png"
$ws.Range("C42").Value = 3
$ws.Range("D42").Value = 0

# Restore view state: zoom level and current selection
# (topLeftCell scroll position is not exposed for writing via this COM host)
$excel.ActiveWindow.ScrollRow = 24
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 139
$ws.Range("B39:E39").Select()
